$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.102.86'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '1.973.66'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  +0.72%  '
$ws.Range("D5").Value = '329.63'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '1.014'
$ws.Range("E6").Value = '  +0.64%  '
$ws.Range("D7").Value = '0.4975'
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("D8").Value = '0.4197'
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").Value = '54.31'
$ws.Range("E9").Value = '  +4.44%  '
$ws.Range("D10").Value = '0.09350'
$ws.Range("E10").Value = '  +4.96%  '
$ws.Range("D11").Value = '1.099'
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("D12").Value = '22.82'
$ws.Range("E12").Value = '  -2.31%  '
$ws.Range("D13").Value = '1.972.81'
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").Value = '7.893'
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("D15").Value = '6.457'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").Value = '1.015'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '0.00001112'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '91.77'
$ws.Range("D19").Value = '0.06709'
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").Value = '19.16'
$ws.Range("E20").Value = '  -3.07%  '
$ws.Range("D21").Value = '1.014'
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").Value = '5.964'
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '29.105.34'
$ws.Range("E23").Value = '  -1.28%  '
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").Value = '2.270'
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").Value = '2.245.90'
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("D27").Value = '20.78'
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").Value = '156.95'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").Value = '6.220'
$ws.Range("E29").Value = '  -4.73%  '
$ws.Range("D30").Value = '2.267'
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("D31").Value = '127.47'
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("D32").Value = '1.047'
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").Value = '0.09840'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").Value = '1.499'
$ws.Range("E34").Value = '  -4.18%  '
$ws.Range("D35").Value = '5.811'
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("D36").Value = '3.749'
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("D37").Value = '0.02418'
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("D38").Value = '1.324'
$ws.Range("E38").Value = '  +3.06%  '
$ws.Range("D39").Value = '0.06405'
$ws.Range("E39").Value = '  +0.97%  '
$ws.Range("D40").Value = '9.057'
$ws.Range("E40").Value = '  -5.47%  '
$ws.Range("D41").Value = '0.6472'
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("D42").Value = '11.51'
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("D43").Value = '0.2004'
$ws.Range("E43").Value = '  -3.13%  '
$ws.Range("D44").Value = '1.013'
$ws.Range("E44").Value = '  +0.71%  '
$ws.Range("D45").Value = '0.6195'
$ws.Range("E45").Value = '  -2.46%  '
$ws.Range("D46").Value = '1.349'
$ws.Range("E46").Value = '  +6.37%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.178'
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '13.23'
$ws.Range("E48").Value = '  -0.98%  '
$ws.Range("D49").Value = '3.484'
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("D50").Value = '0.00000000333'
$ws.Range("E50").Value = '  -0.85%  '
$ws.Range("D51").Value = '0.06959'
$ws.Range("E51").Value = '  -0.42%  '
